$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Get-RowIndexByFirstCellText($table, $text) {
    for ($i = 1; $i -le $table.Rows.Count; $i++) {
        if ($table.Rows.Item($i).Cells.Item(1).Range.Text.TrimEnd([char]7, [char]13) -eq $text) {
            return $i
        }
    }
    return -1
}

# 1) Rename the existing "Class: Anthozoa / corals" row in place to
#    "Class: Anthozoa / jellyfishes" -- oh wait, per the actual target this row becomes
#    "Class: Scyphozoa / jellyfishes" (its 3rd cell "beam and otter trawl" is kept as-is).
$anthozoaIdx = Get-RowIndexByFirstCellText $t "Class: Anthozoa"
$row = $t.Rows.Item($anthozoaIdx)
$row.Cells.Item(1).Range.Text = "Class: Scyphozoa"
$row.Cells.Item(2).Range.Text = "jellyfishes"

# 2) Delete the now-duplicate original "Class: Scyphozoa / jellyfishes / otter trawl" row.
$dupScyphozoaIdx = Get-RowIndexByFirstCellText $t "Class: Scyphozoa"
# There are now two "Class: Scyphozoa" rows -- the renamed one (anthozoaIdx) and the
# original one right after it. Delete the one that is NOT anthozoaIdx, i.e. the next one.
if ($dupScyphozoaIdx -eq $anthozoaIdx) {
    $dupScyphozoaIdx = $anthozoaIdx + 1
}
$t.Rows.Item($dupScyphozoaIdx).Delete()

# 3) Delete the "Family: Pleuronectidae / other flatfishes / otter trawl" row.
$idx = Get-RowIndexByFirstCellText $t "Family: Pleuronectidae"
$t.Rows.Item($idx).Delete()

# 4) Delete the "Class: Asteroidea / other seastars / otter trawl" row.
$idx = Get-RowIndexByFirstCellText $t "Class: Asteroidea"
$t.Rows.Item($idx).Delete()

# 5) Delete the "Class: Holothuroidea / sea cucumbers / beam and otter trawl" row.
$idx = Get-RowIndexByFirstCellText $t "Class: Holothuroidea"
$t.Rows.Item($idx).Delete()

# 6) Add a new row right before "Subphylum: Tunicata" with the renamed Anthozoa entry:
#    "Class: Anthozoa" / "soft corals and sea anemones" / "beam and otter trawl"
$tunicataIdx = Get-RowIndexByFirstCellText $t "Subphylum: Tunicata"
$beforeRow = $t.Rows.Item($tunicataIdx)
$newRow = $t.Rows.Add($beforeRow)

# A freshly Added row only has a single cell; split it twice to get the table's
# 3 columns back.
$newRow.Cells.Item(1).Split(1, 2)
$newRow = $t.Rows.Item($newRow.Index)
$newRow.Cells.Item($newRow.Cells.Count).Split(1, 2)
$newRow = $t.Rows.Item($newRow.Index)

for ($i = 1; $i -le 3; $i++) {
    $cell = $newRow.Cells.Item($i)
    $cell.Range.Style = "Compact"
    $cell.Range.Paragraphs.Item(1).Alignment = 1
}

$newRow.Cells.Item(1).Range.Text = "Class: Anthozoa"
$newRow.Cells.Item(2).Range.Text = "soft corals and sea anemones"
$newRow.Cells.Item(3).Range.Text = "beam and otter trawl"

Write-Host "Final row count: " $t.Rows.Count
